$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "地区名１"
$ws.Range("B2").Value = "地区名２"
$ws.Range("B3").Value = "地区名３"
$ws.Range("B4").Value = "地区名４"
$ws.Range("B5").Value = "地区名５"

$ws.Range("B6").Select()
